# Update the density/depth relationship derived values (peak centers,
# fit errors, deltaNe, Ne_Corr, and residual_1447) for rows 2-13 of the
# Ne line-fitting loop results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1448.004569101116
$ws.Range("D2").Value = 0.00722
$ws.Range("G2").Value = 330.410524708748
$ws.Range("H2").Value = 1.000203039801087
$ws.Range("K2").Value = 141.256915563747

$ws.Range("C3").Value = 1448.001134958084
$ws.Range("D3").Value = 0.00723
$ws.Range("G3").Value = 330.1893144421676
$ws.Range("H3").Value = 1.000873125631825
$ws.Range("K3").Value = 136.6455115309842

$ws.Range("C4").Value = 1448.020926302389
$ws.Range("D4").Value = 0.00702
$ws.Range("G4").Value = 330.1876286963393
$ws.Range("H4").Value = 1.000878235507392
$ws.Range("K4").Value = 133.3476666884557

$ws.Range("C5").Value = 1448.11284983611
$ws.Range("D5").Value = 0.014
$ws.Range("G5").Value = 330.3049332311246
$ws.Range("H5").Value = 1.00052278348721
$ws.Range("K5").Value = 231.059180373549

$ws.Range("C6").Value = 1447.991191605644
$ws.Range("D6").Value = 0.00747
$ws.Range("G6").Value = 330.1878663972166
$ws.Range("H6").Value = 1.000877514979296
$ws.Range("K6").Value = 135.593271084326

$ws.Range("C7").Value = 1448.031301948534
$ws.Range("D7").Value = 0.00663
$ws.Range("G7").Value = 330.190337887884
$ws.Range("H7").Value = 1.000870023362632
$ws.Range("K7").Value = 126.639645308679

$ws.Range("C8").Value = 1448.024610067631
$ws.Range("D8").Value = 0.00707
$ws.Range("G8").Value = 330.199002414817
$ws.Range("H8").Value = 1.00084376021462
$ws.Range("K8").Value = 139.141402826869

$ws.Range("C9").Value = 1448.034876468905
$ws.Range("D9").Value = 0.00645
$ws.Range("G9").Value = 330.1887387324909
$ws.Range("H9").Value = 1.000874870731751
$ws.Range("K9").Value = 121.9973720959166

$ws.Range("C10").Value = 1448.029871734082
$ws.Range("D10").Value = 0.00676
$ws.Range("G10").Value = 330.1802150671983
$ws.Range("H10").Value = 1.000900708507634
$ws.Range("K10").Value = 128.0108410899661

$ws.Range("C11").Value = 1448.024917666296
$ws.Range("D11").Value = 0.00708
$ws.Range("G11").Value = 330.1896489436203
$ws.Range("H11").Value = 1.00087211168873
$ws.Range("K11").Value = 134.6010612846933

$ws.Range("C12").Value = 1448.022255336594
$ws.Range("D12").Value = 0.00673
$ws.Range("G12").Value = 330.1895302861935
$ws.Range("H12").Value = 1.0008724713637
$ws.Range("K12").Value = 125.9009063268143

$ws.Range("C13").Value = 1448.017430005813
$ws.Range("D13").Value = 0.0071
$ws.Range("G13").Value = 330.1849733290023
$ws.Range("H13").Value = 1.000886284630124
$ws.Range("K13").Value = 135.9233781216712
